$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts old A -> B, old B -> C)
$ws.Columns("A").Insert()

# Header row
$ws.Range("A1").Value = "Дата и время измерения"

# Dates / timestamps for existing rows (as text)
$ws.Range("A2").Value = "2020-04-03 02:10"
$ws.Range("A3").Value = "2020-04-03 02:10"
$ws.Range("A4").Value = "2020-04-03 03:47:39"
$ws.Range("A5").Value = "2020-04-03 06:47:07"
$ws.Range("A6").Value = "2020-04-10 09:44:47"
$ws.Range("A7").Value = "2020-04-10 10:18:12"
$ws.Range("A8").Value = "2020-04-10 13:18:07"
$ws.Range("A9").Value = "2020-04-10 13:18:11"
$ws.Range("A10").Value = "2020-04-10 13:20:24"
$ws.Range("A11").Value = "2020-04-10 13:20:31"

# New rows appended at the bottom
$ws.Range("A12").Value = "2020-04-10 13:30:42"
$ws.Range("B12").Value = 177
$ws.Range("C12").Value = 77

$ws.Range("A13").Value = "2020-04-10 13:43:14"
$ws.Range("B13").Value = 110
$ws.Range("C13").Value = 70

$ws.Range("A14").Value = "2020-04-14 00:22:19"
$ws.Range("B14").Value = 114
$ws.Range("C14").Value = 80
